$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '322.99'
$c.NumberFormat = "General"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '-2.76%'
$c.NumberFormat = "General"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '42.60'
$c.NumberFormat = "General"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '-6.87%'
$c.NumberFormat = "General"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '5.267'
$c.NumberFormat = "General"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '-7.51%'
$c.NumberFormat = "General"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.08136'
$c.NumberFormat = "General"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '-2.80%'
$c.NumberFormat = "General"

$ws.Range("B6").Value = 'FTXToken'

$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.798'
$c.NumberFormat = "General"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '-11.77%'
$c.NumberFormat = "General"

$ws.Range("B7").Value = 'MXToken'

$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.9568'
$c.NumberFormat = "General"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '-2.49%'
$c.NumberFormat = "General"

$ws.Range("B8").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C8").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.1127'
$c.NumberFormat = "General"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '-3.39%'
$c.NumberFormat = "General"

$ws.Range("B9").Value = 'WazirX'

$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.1864'
$c.NumberFormat = "General"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '-4.05%'
$c.NumberFormat = "General"

$ws.Range("B10").Value = 'MandalaExchangeToken'

$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.09323'
$c.NumberFormat = "General"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '-7.63%'
$c.NumberFormat = "General"

$ws.Range("B11").Value = 'BitrueCoin'

$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.04619'
$c.NumberFormat = "General"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '-1.04%'
$c.NumberFormat = "General"

$ws.Range("B12").Value = 'MCDex'

$ws.Range("C12").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.507'
$c.NumberFormat = "General"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '-27.55%'
$c.NumberFormat = "General"

$ws.Range("B13").Value = 'BitMartToken'

$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.1062'
$c.NumberFormat = "General"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '0.40%'
$c.NumberFormat = "General"

$ws.Range("B14").Value = 'BitForexToken'

$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.001288'
$c.NumberFormat = "General"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '-0.09%'
$c.NumberFormat = "General"

$ws.Range("B15").Value = 'TigerCash'

$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.005810'
$c.NumberFormat = "General"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '-4.08%'
$c.NumberFormat = "General"

$ws.Range("B16").Value = 'LEO'

$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '3.378'
$c.NumberFormat = "General"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '0.24%'
$c.NumberFormat = "General"

$ws.Range("B17").Value = 'GateToken'

$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '4.377'
$c.NumberFormat = "General"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '-2.05%'
$c.NumberFormat = "General"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.540'
$c.NumberFormat = "General"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '-2.19%'
$c.NumberFormat = "General"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.3364'
$c.NumberFormat = "General"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '0.40%'
$c.NumberFormat = "General"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.1392'
$c.NumberFormat = "General"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '-0.58%'
$c.NumberFormat = "General"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.2596'
$c.NumberFormat = "General"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '0.11%'
$c.NumberFormat = "General"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.04209'
$c.NumberFormat = "General"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '-0.02%'
$c.NumberFormat = "General"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.001264'
$c.NumberFormat = "General"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '-3.33%'
$c.NumberFormat = "General"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.004310'
$c.NumberFormat = "General"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '-6.08%'
$c.NumberFormat = "General"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.0001312'
$c.NumberFormat = "General"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '2.43%'
$c.NumberFormat = "General"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.0003006'
$c.NumberFormat = "General"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '-19.66%'
$c.NumberFormat = "General"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.02596'
$c.NumberFormat = "General"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '-6.56%'
$c.NumberFormat = "General"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.05462'
$c.NumberFormat = "General"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '-6.13%'
$c.NumberFormat = "General"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.007831'
$c.NumberFormat = "General"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '1.19%'
$c.NumberFormat = "General"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.1394'
$c.NumberFormat = "General"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '-3.09%'
$c.NumberFormat = "General"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.006606'
$c.NumberFormat = "General"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '-8.18%'
$c.NumberFormat = "General"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.002133'
$c.NumberFormat = "General"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '8.05%'
$c.NumberFormat = "General"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.008695'
$c.NumberFormat = "General"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '6.29%'
$c.NumberFormat = "General"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.3300'
$c.NumberFormat = "General"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.00007017'
$c.NumberFormat = "General"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '-2.60%'
$c.NumberFormat = "General"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.00000000757'
$c.NumberFormat = "General"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '0.87%'
$c.NumberFormat = "General"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.003498'
$c.NumberFormat = "General"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '0.15%'
$c.NumberFormat = "General"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.003562'
$c.NumberFormat = "General"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '1.79%'
$c.NumberFormat = "General"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.00002119'
$c.NumberFormat = "General"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '0.87%'
$c.NumberFormat = "General"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0002018'
$c.NumberFormat = "General"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '0.87%'
$c.NumberFormat = "General"
